$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149 (shifts existing rows 149:162 down to 150:163,
# matching the diff's net effect of one new weekly record at the top of the
# "Vega Modelo de Temuco - Pomelo" block).
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new weekly record.
$ws.Cells.Item(149, 1).Value  = 10
$ws.Cells.Item(149, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(149, 3).Value  = 'La Araucanía'
$ws.Cells.Item(149, 4).Value  = 44505
$ws.Cells.Item(149, 5).Value  = 9
$ws.Cells.Item(149, 6).Value  = 'Fruta'
$ws.Cells.Item(149, 7).Value  = 100102
$ws.Cells.Item(149, 8).Value  = 'Cítricos'
$ws.Cells.Item(149, 9).Value  = 100102006
$ws.Cells.Item(149, 10).Value = 'Pomelo'
$ws.Cells.Item(149, 11).Value = 'Start Ruby'
$ws.Cells.Item(149, 12).Value = 'Primera'
$ws.Cells.Item(149, 13).Value = 75
$ws.Cells.Item(149, 14).Value = 12000
$ws.Cells.Item(149, 15).Value = 13000
$ws.Cells.Item(149, 16).Value = 12467
$ws.Cells.Item(149, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(149, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(149, 19).Value = 831
$ws.Cells.Item(149, 20).Value = 15
